$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue([string]$cellRef, [string]$val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

$ws.Range("D2").Value = "25.309.76"
$ws.Range("E2").Value = "  -2.73%  "
$ws.Range("D3").Value = "1.571.96"
$ws.Range("E3").Value = "  -3.79%  "
$ws.Range("E4").Value = "  -0.12%  "
Set-TextValue "D5" "208.09"
$ws.Range("E5").Value = "  -3.08%  "
Set-TextValue "D7" "0.478"
$ws.Range("E7").Value = "  -4.84%  "
$ws.Range("E8").Value = "  -2.40%  "
$ws.Range("E9").Value = "  -1.92%  "
Set-TextValue "D10" "17.99"
$ws.Range("E10").Value = "  -2.14%  "
Set-TextValue "D11" "0.0781"
$ws.Range("E11").Value = "  -1.30%  "
$ws.Range("D12").Value = "1.790.54"
$ws.Range("E12").Value = "  -3.76%  "
$ws.Range("D13").Value = "1.574.65"
$ws.Range("E13").Value = "  -3.62%  "
$ws.Range("E14").Value = "  -3.33%  "
$ws.Range("E15").Value = "  -3.84%  "
$ws.Range("D16").Value = "25.323.59"
$ws.Range("E16").Value = "  -2.62%  "
Set-TextValue "D17" "59.78"
$ws.Range("E17").Value = "  -2.70%  "
$ws.Range("D18").Value = "0.0₃0710"
$ws.Range("E18").Value = "  -4.20%  "
$ws.Range("E19").Value = "  -0.07%  "
Set-TextValue "D20" "185.44"
$ws.Range("E20").Value = "  -2.82%  "
Set-TextValue "D21" "4.14"
$ws.Range("E21").Value = "  -2.46%  "
$ws.Range("E22").Value = "  -3.13%  "
$ws.Range("E23").Value = "  -3.19%  "
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("E25").Value = "  -2.57%  "
Set-TextValue "D26" "141.21"
$ws.Range("E26").Value = "  -1.69%  "
$ws.Range("E27").Value = "  -5.30%  "
$ws.Range("E28").Value = "  -4.47%  "
Set-TextValue "D29" "14.89"
$ws.Range("E29").Value = "  -1.89%  "
$ws.Range("E30").Value = "  -6.41%  "
$ws.Range("E31").Value = "  -4.28%  "
$ws.Range("E32").Value = "  -2.59%  "
Set-TextValue "D33" "3.02"
$ws.Range("E33").Value = "  -3.04%  "
$ws.Range("E34").Value = "  -1.44%  "
Set-TextValue "D35" "2.25"
$ws.Range("E35").Value = "  -6.31%  "
$ws.Range("D36").Value = "1.087.90"
$ws.Range("E36").Value = "  -3.89%  "
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D38" "0.0151"
$ws.Range("E38").Value = "  -2.47%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D39" "2.32"
$ws.Range("E39").Value = "  -4.83%  "
$ws.Range("E40").Value = "  -9.50%  "
Set-TextValue "D41" "0.495"
$ws.Range("E41").Value = "  -4.37%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D42" "0.758"
$ws.Range("E42").Value = "  -2.22%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D43" "93.44"
$ws.Range("E43").Value = "  -5.03%  "
$ws.Range("E44").Value = "  -3.00%  "
$ws.Range("D45").Value = "1.704.42"
$ws.Range("E45").Value = "  -3.75%  "
$ws.Range("D46").Value = "0.0₆0112"
$ws.Range("E46").Value = "  -2.41%  "
Set-TextValue "D47" "52.81"
$ws.Range("E47").Value = "  -3.66%  "
$ws.Range("E48").Value = "  -3.64%  "
$ws.Range("E49").Value = "  -4.35%  "
$ws.Range("E50").Value = "  -1.73%  "
$ws.Range("E51").Value = "  -0.26%  "
